# Update "想去人数" (want-to-go count) figures on the 展览, 演出 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 754
$ws1.Range("F3").Value = 669
$ws1.Range("F4").Value = 25
$ws1.Range("F6").Value = 1184
$ws1.Range("F8").Value = 32
$ws1.Range("F9").Value = 46
$ws1.Range("F10").Value = 591
$ws1.Range("F14").Value = 98
$ws1.Range("F16").Value = 88
$ws1.Range("F17").Value = 292
$ws1.Range("F18").Value = 400
$ws1.Range("F19").Value = 489
$ws1.Range("F21").Value = 6101
$ws1.Range("F22").Value = 5303

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 90

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 754
$ws4.Range("F3").Value = 669
$ws4.Range("F4").Value = 25
$ws4.Range("F6").Value = 1184
$ws4.Range("F8").Value = 32
$ws4.Range("F9").Value = 46
$ws4.Range("F10").Value = 591
$ws4.Range("F13").Value = 90
$ws4.Range("F16").Value = 98
$ws4.Range("F18").Value = 88
$ws4.Range("F19").Value = 292
$ws4.Range("F20").Value = 400
$ws4.Range("F21").Value = 489
$ws4.Range("F23").Value = 6101
$ws4.Range("F25").Value = 5303
